# Auto-applied data refresh: update per-leve market-price / profit
# figures (columns H-N) across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets,
# mirroring the scheduled Universalis price-refresh job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 6
$ws.Range("H6").Value = 389103.16
$ws.Range("I6").Value = 626329.9
$ws.Range("J6").Value = 9540.4
$ws.Range("K6").Value = 1878989.7
$ws.Range("L6").Value = 28621.2
$ws.Range("M6").Value = -1878877.7
$ws.Range("N6").Value = -28845.2

# Row 64
$ws.Range("H64").Value = 27053.334
$ws.Range("I64").Value = 3345
$ws.Range("J64").Value = 102920
$ws.Range("K64").Value = 3345
$ws.Range("L64").Value = 102920
$ws.Range("M64").Value = -3097
$ws.Range("N64").Value = -103416

# Row 67
$ws.Range("H67").Value = 27053.334
$ws.Range("I67").Value = 3345
$ws.Range("J67").Value = 102920
$ws.Range("K67").Value = 3345
$ws.Range("L67").Value = 102920
$ws.Range("M67").Value = -2487
$ws.Range("N67").Value = -104636

# Row 76
$ws.Range("H76").Value = 100003720
$ws.Range("I76").Value = 142860180
$ws.Range("J76").Value = 5334.6665
$ws.Range("K76").Value = 142860180
$ws.Range("L76").Value = 5334.6665
$ws.Range("M76").Value = -142859865
$ws.Range("N76").Value = -5964.6665

# Row 79
$ws.Range("H79").Value = 100003720
$ws.Range("I79").Value = 142860180
$ws.Range("J79").Value = 5334.6665
$ws.Range("K79").Value = 142860180
$ws.Range("L79").Value = 5334.6665
$ws.Range("M79").Value = -142859088
$ws.Range("N79").Value = -7518.6665

# Row 96
$ws.Range("H96").Value = 604.1429000000001
$ws.Range("I96").Value = 650
$ws.Range("J96").Value = 329
$ws.Range("K96").Value = 1950
$ws.Range("L96").Value = 987
$ws.Range("M96").Value = -577
$ws.Range("N96").Value = -3733

# Row 100
$ws.Range("H100").Value = 5267.5557
$ws.Range("I100").Value = 3762.3076
$ws.Range("J100").Value = 9181.200000000001
$ws.Range("K100").Value = 3762.3076
$ws.Range("L100").Value = 9181.200000000001
$ws.Range("M100").Value = -3221.3076
$ws.Range("N100").Value = -10263.2

# Row 118
$ws.Range("H118").Value = 1042.32
$ws.Range("I118").Value = 857.7778
$ws.Range("J118").Value = 1146.125
$ws.Range("K118").Value = 2573.3334
$ws.Range("L118").Value = 3438.375
$ws.Range("M118").Value = -916.3334
$ws.Range("N118").Value = -6752.375

$ws = $wb.Worksheets("ARM")
# Row 32
$ws.Range("H32").Value = 19594.064
$ws.Range("I32").Value = 3296.3584
$ws.Range("J32").Value = 115569.445
$ws.Range("K32").Value = 3296.3584
$ws.Range("L32").Value = 115569.445
$ws.Range("M32").Value = -3009.3584
$ws.Range("N32").Value = -116143.445

# Row 46
$ws.Range("H46").Value = 3660
$ws.Range("I46").Value = 2800
$ws.Range("J46").Value = 3875
$ws.Range("K46").Value = 2800
$ws.Range("L46").Value = 3875
$ws.Range("M46").Value = -2481
$ws.Range("N46").Value = -4513

# Row 61
$ws.Range("H61").Value = 1529.3334
$ws.Range("I61").Value = 1505.6364
$ws.Range("J61").Value = 1790
$ws.Range("K61").Value = 1505.6364
$ws.Range("L61").Value = 1790
$ws.Range("M61").Value = -1293.6364
$ws.Range("N61").Value = -2214

# Row 74
$ws.Range("H74").Value = 48590.906
$ws.Range("I74").Value = 50995.95
$ws.Range("J74").Value = 490
$ws.Range("K74").Value = 50995.95
$ws.Range("L74").Value = 490
$ws.Range("M74").Value = -50121.95
$ws.Range("N74").Value = -2238

# Row 77
$ws.Range("H77").Value = 48590.906
$ws.Range("I77").Value = 50995.95
$ws.Range("J77").Value = 490
$ws.Range("K77").Value = 254979.75
$ws.Range("L77").Value = 2450
$ws.Range("M77").Value = -250611.75
$ws.Range("N77").Value = -11186

# Row 136
$ws.Range("H136").Value = 1529.3334
$ws.Range("I136").Value = 1505.6364
$ws.Range("J136").Value = 1790
$ws.Range("K136").Value = 4516.9092
$ws.Range("L136").Value = 5370
$ws.Range("M136").Value = -1966.9092
$ws.Range("N136").Value = -10470

$ws = $wb.Worksheets("BSM")
# Row 88
$ws.Range("H88").Value = 39233
$ws.Range("I88").Value = 9749.5
$ws.Range("J88").Value = 58888.668
$ws.Range("K88").Value = 9749.5
$ws.Range("L88").Value = 58888.668
$ws.Range("M88").Value = -9343.5
$ws.Range("N88").Value = -59700.668

# Row 91
$ws.Range("H91").Value = 39233
$ws.Range("I91").Value = 9749.5
$ws.Range("J91").Value = 58888.668
$ws.Range("K91").Value = 9749.5
$ws.Range("L91").Value = 58888.668
$ws.Range("M91").Value = -8345.5
$ws.Range("N91").Value = -61696.668

$ws = $wb.Worksheets("CUL")
# Row 5
$ws.Range("H5").Value = 6582.4116
$ws.Range("I5").Value = 450.16666
$ws.Range("K5").Value = 1350.49998
$ws.Range("M5").Value = -1238.49998

# Row 122
$ws.Range("H122").Value = 644.1429000000001
$ws.Range("J122").Value = 1065
$ws.Range("L122").Value = 9585
$ws.Range("N122").Value = -14485

# Row 135
$ws.Range("H135").Value = 6582.4116
$ws.Range("I135").Value = 450.16666
$ws.Range("K135").Value = 4051.49994
$ws.Range("M135").Value = -1516.49994

$ws = $wb.Worksheets("GSM")
# Row 70
$ws.Range("H70").Value = 4146.28
$ws.Range("J70").Value = 4248.625
$ws.Range("L70").Value = 4248.625
$ws.Range("N70").Value = -4788.625

# Row 73
$ws.Range("H73").Value = 4146.28
$ws.Range("J73").Value = 4248.625
$ws.Range("L73").Value = 4248.625
$ws.Range("N73").Value = -6120.625

$ws = $wb.Worksheets("LTW")
# Row 136
$ws.Range("H136").Value = 3558
$ws.Range("I136").Value = 3874.5
$ws.Range("J136").Value = 3377.1428
$ws.Range("K136").Value = 11623.5
$ws.Range("L136").Value = 10131.4284
$ws.Range("M136").Value = -9073.5
$ws.Range("N136").Value = -15231.4284

$ws = $wb.Worksheets("WVR")
# Row 63
$ws.Range("H63").Value = 23333.334
$ws.Range("J63").Value = 23333.334
$ws.Range("L63").Value = 23333.334
$ws.Range("N63").Value = -24581.334

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 66
$ws.Range("H66").Value = 23333.334
$ws.Range("J66").Value = 23333.334
$ws.Range("L66").Value = 70000.00199999999
$ws.Range("N66").Value = -76240.00199999999

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 69
$ws.Range("H69").Value = 10125
$ws.Range("J69").Value = 10125
$ws.Range("L69").Value = 10125
$ws.Range("N69").Value = -11623

# Row 70
$ws.Range("H70").Value = 33105
$ws.Range("J70").Value = 33105
$ws.Range("L70").Value = 33105
$ws.Range("N70").Value = -33735

# Row 72
$ws.Range("H72").Value = 10125
$ws.Range("J72").Value = 10125
$ws.Range("L72").Value = 30375
$ws.Range("N72").Value = -37863

# Row 73
$ws.Range("H73").Value = 33105
$ws.Range("J73").Value = 33105
$ws.Range("L73").Value = 33105
$ws.Range("N73").Value = -35289

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

# Row 76
$ws.Range("H76").Value = 17750
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 17750
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 17750
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -18380

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

# Row 79
$ws.Range("H79").Value = 17750
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 17750
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 17750
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -19934

# Row 82
$ws.Range("H82").Value = 40075.25
$ws.Range("J82").Value = 40075.25
$ws.Range("L82").Value = 40075.25
$ws.Range("N82").Value = -40841.25

# Row 85
$ws.Range("H85").Value = 40075.25
$ws.Range("J85").Value = 40075.25
$ws.Range("L85").Value = 40075.25
$ws.Range("N85").Value = -42727.25

# Row 136
$ws.Range("H136").Value = 2978246
$ws.Range("I136").Value = 2978246
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8934738
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8932188
$ws.Range("N136").ClearContents()
